$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.033.71"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.643.94"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").Value = "'216.51"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.0640"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.870.79"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "1.649.81"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "0.0₃0768"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'63.11"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "25.948.94"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "'9.95"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "'144.86"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "'6.95"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("D36").Value = "'0.905"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "1.134.90"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'99.71"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "1.779.63"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("D46").Value = "'56.86"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "'1.47"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +0.00%  "
